$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: merge the two runs "WED Oct 04" + " 13:43:52 PDT 2017"
# into a single run "WED Oct 04 13:43:52 PDT 2017".
# Executing Find/Replace across the two runs causes Word to coalesce
# the matched text into one run using the first run's formatting.
# ------------------------------------------------------------------
$d.Content.Find.Execute("WED Oct 04 13:43:52 PDT 2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "WED Oct 04 13:43:52 PDT 2017", 2)

# ------------------------------------------------------------------
# Change 2: append a new transaction block (THU Oct 05 ...) right
# after the existing "Amount balance" paragraph, before the two
# trailing blank paragraphs at the end of the document.
# ------------------------------------------------------------------
$amountBalancePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Amount balance*") {
        $amountBalancePara = $p
    }
}
$insertPoint = $d.Range($amountBalancePara.Range.End, $amountBalancePara.Range.End)

$newBlockXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>THU Oct 05</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t xml:space="preserve"> 11:02:51 PDT 2017</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>Person Name</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/><w:t>- NS M</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>---------------------------------------------------------------</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>Item Name</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/><w:t>- CARROT</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/></w:rPr><w:t>Amount Received</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/></w:rPr><w:tab/><w:t>- 2232</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>Amount Received mode</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:tab/><w:t>- CASH AND CLEARD</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PlainText"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($newBlockXml)
